$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 263, shifting existing rows 263..331 down to 264..332.
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row 263 with the new record
# (Vega Monumental Concepción - Piña, Caramelo, Segunda).
$ws.Range("A263").Value = 11
$ws.Range("B263").Value = "Vega Monumental Concepción"
$ws.Range("C263").Value = "Bíobío"
$ws.Range("D263").Value = 45275
$ws.Range("E263").Value = 8
$ws.Range("F263").Value = "Fruta"
$ws.Range("G263").Value = 100108
$ws.Range("H263").Value = "Tropicales y subtropicales"
$ws.Range("I263").Value = 100108005
$ws.Range("J263").Value = "Piña"
$ws.Range("K263").Value = "Caramelo"
$ws.Range("L263").Value = "Segunda"
$ws.Range("M263").Value = 200
$ws.Range("N263").Value = 24000
$ws.Range("O263").Value = 25000
$ws.Range("P263").Value = 24500
$ws.Range("Q263").Value = "`$/caja 14 unidades"
$ws.Range("R263").Value = "Ecuador"
$ws.Range("S263").Value = 1750
$ws.Range("T263").Value = 14
